$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before the existing row 81 ("li: Roberts McCubbin Primary
#    School, BOX HILL SOUTH") with the new Mooroolbark College entry. This
#    shifts every subsequent row down by one.
$ws.Rows.Item(81).Insert()
$ws.Range("A81").Value = "li: Mooroolbark College, MOOROOLBARK"

# 2. Update the big concatenated school-closures cell (now at row 109 after the
#    insert above) to add "Lalor Secondary College, LALOR" between Fitzroy
#    Primary School and Macleod College.
$ws.Range("A109").Value = "Al Siraat College, EPPINGCharles La Trobe, MACLEOD WESTDiamond Valley College, DIAMOND REEKEpping Secondary College, EPPINGFitzroy High School, FITZROYFitzroy Primary School, FITZROYLalor Secondary College, LALORMacleod College, MACLEODNewbury Primary School, CRAGIEBURNNorthcote High School, NORTHCOTENorthern College of Arts and Technology (NCAT), PRESTONPascoe Vale Girls Secondary College, OAK PARKPenders Grove Primary School, THORNBURYPrinces Hill Secondary College, PRINCES HILLSacred Heart School, FITZROYRoxburgh College, ROXBURGH PARKWatsonia Primary School, WATSONIATAFE"

# 3. Insert a new row before the existing row 170 ("li: Parkville College
#    (Parkville and Malmsbury Campus)") with the new Mary McKillop College
#    entry. This shifts every subsequent row down by one again.
$ws.Rows.Item(170).Insert()
$ws.Range("A170").Value = "li: Mary McKillop College, WERRIBEE"

# 4. Update the "Last Update" row (now at row 188 after both inserts) with the
#    new date.
$ws.Range("A188").Value = "li: Last Update: 26 July 2020"
